$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value, $fmtSourceCellRef) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
    $ws.Range($fmtSourceCellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

$ws.Range("E2").Value = "2026-02-25 21:18:39"
$ws.Range("O2").Value = "5.5 °C"
$ws.Range("E3").Value = "2026-02-25 21:18:41"
$ws.Range("E4").Value = "2026-02-25 21:18:44"
$ws.Range("J4").Value = "1021.9 hPa"
$ws.Range("E5").Value = "2026-02-25 21:18:46"
$ws.Range("N5").Value = "2.0 °C 20:54 TU"
$ws.Range("O5").Value = "5.7 °C"
$ws.Range("E6").Value = "2026-02-25 21:18:49"
$ws.Range("J6").Value = "1021.8 hPa"
$ws.Range("E7").Value = "2026-02-25 21:18:52"
$ws.Range("J7").Value = "1021.4 hPa"
$ws.Range("E8").Value = "2026-02-25 21:18:54"
$ws.Range("J8").Value = "1021.1 hPa"
$ws.Range("E9").Value = "2026-02-25 21:18:57"
$ws.Range("E10").Value = "2026-02-25 21:19:00"
$ws.Range("O10").Value = "9.4 °C"
$ws.Range("E11").Value = "2026-02-25 21:19:02"
$ws.Range("O11").Value = "8.9 °C"
$ws.Range("E12").Value = "2026-02-25 21:19:05"
$ws.Range("E13").Value = "2026-02-25 21:19:07"
$ws.Range("J13").Value = "1022.8 hPa"
$ws.Range("O13").Value = "6.8 °C"
$ws.Range("E14").Value = "2026-02-25 21:19:10"
$ws.Range("E15").Value = "2026-02-25 21:19:12"
$ws.Range("O15").Value = "10.1 °C"
$ws.Range("E16").Value = "2026-02-25 21:19:14"
$ws.Range("N16").Value = "0.7 °C 20:48 TU"
$ws.Range("E17").Value = "2026-02-25 21:19:17"
$ws.Range("O17").Value = "9.2 °C"
$ws.Range("E18").Value = "2026-02-25 21:19:20"
$ws.Range("J18").Value = "1022.0 hPa"
$ws.Range("E19").Value = "2026-02-25 21:19:22"
$ws.Range("N19").Value = "8.0 °C 20:52 TU"
$ws.Range("O19").Value = "12.2 °C"
$ws.Range("E20").Value = "2026-02-25 21:19:25"
$ws.Range("K20").Value = "16.6 MJ/m2"
$ws.Range("N20").Value = "-0.3 °C 20:59 TU"
$ws.Range("E21").Value = "2026-02-25 21:19:28"
$ws.Range("J21").Value = "1021.5 hPa"
$ws.Range("E22").Value = "2026-02-25 21:19:30"
$ws.Range("N22").Value = "0.1 °C 20:57 TU"
$ws.Range("E23").Value = "2026-02-25 21:19:33"
$ws.Range("N23").Value = "1.5 °C 20:59 TU"
$ws.Range("E24").Value = "2026-02-25 21:19:35"
$ws.Range("J24").Value = "1020.1 hPa"
$ws.Range("K24").Value = "15.2 MJ/m2"
$ws.Range("O24").Value = "10.8 °C"
$ws.Range("E25").Value = "2026-02-25 21:19:38"
$ws.Range("O25").Value = "5.2 °C"
$ws.Range("E26").Value = "2026-02-25 21:19:41"
$ws.Range("J26").Value = "1019.6 hPa"
$ws.Range("N26").Value = "5.3 °C 20:59 TU"
$ws.Range("O26").Value = "10.0 °C"
$ws.Range("E27").Value = "2026-02-25 21:19:43"
$ws.Range("N27").Value = "2.6 °C 20:50 TU"
$ws.Range("E28").Value = "2026-02-25 21:19:46"
$ws.Range("J28").Value = "1021.9 hPa"
$ws.Range("E29").Value = "2026-02-25 21:19:49"
$ws.Range("E30").Value = "2026-02-25 21:19:51"
$ws.Range("J30").Value = "1021.9 hPa"
$ws.Range("E31").Value = "2026-02-25 21:19:54"
$ws.Range("J31").Value = "1021.5 hPa"
$ws.Range("E32").Value = "2026-02-25 21:19:56"
$ws.Range("K32").Value = "16.5 MJ/m2"
$ws.Range("O32").Value = "9.3 °C"
$ws.Range("E33").Value = "2026-02-25 21:19:59"
$ws.Range("J33").Value = "1021.3 hPa"
$ws.Range("K33").Value = "15.5 MJ/m2"
$ws.Range("O33").Value = "8.5 °C"
$ws.Range("E34").Value = "2026-02-25 21:20:02"
$ws.Range("O34").Value = "3.4 °C"
$ws.Range("E35").Value = "2026-02-25 21:20:04"
$ws.Range("J35").Value = "1019.6 hPa"
$ws.Range("E36").Value = "2026-02-25 21:20:07"
$ws.Range("J36").Value = "1022.0 hPa"
$ws.Range("O36").Value = "11.3 °C"
$ws.Range("E37").Value = "2026-02-25 21:20:10"
$ws.Range("J37").Value = "1023.5 hPa"
$ws.Range("O37").Value = "6.6 °C"
$ws.Range("E38").Value = "2026-02-25 21:20:12"
$ws.Range("E39").Value = "2026-02-25 21:20:15"
$ws.Range("O39").Value = "2.5 °C"
$ws.Range("E40").Value = "2026-02-25 21:20:17"
$ws.Range("J40").Value = "1021.9 hPa"
$ws.Range("O40").Value = "9.6 °C"
$ws.Range("E41").Value = "2026-02-25 21:20:20"
$ws.Range("E42").Value = "2026-02-25 21:20:23"
$ws.Range("E43").Value = "2026-02-25 21:20:25"
$ws.Range("E44").Value = "2026-02-25 21:20:28"
$ws.Range("E45").Value = "2026-02-25 21:20:31"
$ws.Range("J45").Value = "1019.9 hPa"
$ws.Range("O45").Value = "10.9 °C"
$ws.Range("E46").Value = "2026-02-25 21:20:33"
$ws.Range("J46").Value = "1020.8 hPa"

Set-TextValue "H2" "47%" "G2"
Set-TextValue "H5" "29%" "G5"
Set-TextValue "H12" "97%" "G12"
Set-TextValue "H14" "91%" "G14"
Set-TextValue "H16" "32%" "G16"
Set-TextValue "H20" "50%" "G20"
Set-TextValue "H22" "44%" "G22"
Set-TextValue "H32" "51%" "G32"
Set-TextValue "H34" "53%" "G34"
Set-TextValue "H43" "71%" "G43"
Set-TextValue "H45" "41%" "G45"
